$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("TestData")

# 1. Row 8: keep A8 ("Search Query") but change B8 from the text "Selenium"
#    to the numeric value 12, and tighten the row height to 13.8 (matches
#    the style already used by row 3).
$ws.Range("B8").Value = 12
$ws.Rows.Item(8).RowHeight = 13.8

# 2. Insert a new row 9 with A9 = "Search Query 2" (renamed from the old
#    "Selenium" shared string) and B9 = 12.2.
$ws.Range("A9").Value = "Search Query 2"
$ws.Range("B9").Value = 12.2
$ws.Rows.Item(9).RowHeight = 15

# 3. Update the active selection to B10, as recorded in the saved view.
[void]$ws.Range("B10").Select()
